$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting everything down by one.
# After the insert, the former grand-total row (old row 16, "=SUM(E2:E15)")
# is now row 17, and Excel has already adjusted its formula to
# "=SUM(E3:E16)" to keep pointing at the same (shifted) data block.
$ws.Rows("1:1").Insert()

# Move the grand-total formula out of row 17 up into the new row 1 (E1),
# leaving row 17's E cell blank but still carrying its original
# (bold/coloured "total") style.
$ws.Range("E1").Formula = $ws.Range("E17").Formula
$ws.Range("E17").Formula = $null

# E1 uses the same currency number format (164) as the totals elsewhere,
# but with the plain/default font, fill and border - not the bold
# highlighted look of the in-table total cell.
$ws.Range("E1").NumberFormat = "#,##0.00\ " + [char]34 + "€" + [char]34
$ws.Range("E1").Font.Bold = $false

# Re-freeze the header/title rows at the new split point (row 2, i.e. below
# the title row and the header row) and restore the selection/scroll
# position, shifted down by one row same as everything else.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E23").Select() | Out-Null

# The AutoFilter range and the (stale, cached) sort-state range both need to
# move down by one row too.
$ws.AutoFilterMode = $false
$ws.Range("A2:Q494").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name Excel maintains alongside the
# AutoFilter also needs to be repointed at the shifted range.
foreach ($n in $ws.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Vorgangsuebersicht!`$A`$2:`$Q`$494"
    }
}

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A2499")) | Out-Null
$ws.Sort.SetRange($ws.Range("A3:Q2499"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

